$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 1389.2858  # H19: 1557.2858 -> 1389.2858
$ws.Cells.Item(19, 9).Value = 1516.6666  # I19: 1637.75 -> 1516.6666
$ws.Cells.Item(19, 10).Value = 1293.75  # J19: 1450 -> 1293.75
$ws.Cells.Item(19, 11).Value = 1516.6666  # K19: 1637.75 -> 1516.6666
$ws.Cells.Item(19, 12).Value = 1293.75  # L19: 1450 -> 1293.75
$ws.Cells.Item(19, 13).Value = -1341.6666  # M19: -1462.75 -> -1341.6666
$ws.Cells.Item(19, 14).Value = -1643.75  # N19: -1800 -> -1643.75
$ws.Cells.Item(32, 8).Value = 799.6667  # H32: 800 -> 799.6667
$ws.Cells.Item(32, 9).Value = 799  # I32: 0 -> 799
$ws.Cells.Item(32, 11).Value = 799  # K32: 0 -> 799
$ws.Cells.Item(32, 13).Value = -473  # M32: None -> -473
$ws.Cells.Item(127, 8).Value = 1090  # H127: 1633.3334 -> 1090
$ws.Cells.Item(127, 9).Value = 1237.5  # I127: 1633.3334 -> 1237.5
$ws.Cells.Item(127, 10).Value = 500  # J127: 0 -> 500
$ws.Cells.Item(127, 11).Value = 3712.5  # K127: 4900.0002 -> 3712.5
$ws.Cells.Item(127, 12).Value = 1500  # L127: 0 -> 1500
$ws.Cells.Item(127, 13).Value = 1247.5  # M127: 59.9997999999996 -> 1247.5
$ws.Cells.Item(127, 14).Value = -11420  # N127: None -> -11420
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2902.9092  # H45: 2745.3333 -> 2902.9092
$ws.Cells.Item(45, 9).Value = 1943  # I45: 1826.625 -> 1943
$ws.Cells.Item(45, 11).Value = 1943  # K45: 1826.625 -> 1943
$ws.Cells.Item(45, 13).Value = -1566  # M45: -1449.625 -> -1566
$ws.Cells.Item(61, 8).Value = 3757.6667  # H61: 4193.375 -> 3757.6667
$ws.Cells.Item(61, 9).Value = 3757.6667  # I61: 4309.4 -> 3757.6667
$ws.Cells.Item(61, 10).Value = 0  # J61: 4000 -> 0
$ws.Cells.Item(61, 11).Value = 3757.6667  # K61: 4309.4 -> 3757.6667
$ws.Cells.Item(61, 12).Value = 0  # L61: 4000 -> 0
$ws.Cells.Item(61, 13).Value = -3545.6667  # M61: -4097.4 -> -3545.6667
$ws.Cells.Item(61, 14).Value = ""  # N61: -4424 -> (blank)
$ws.Cells.Item(136, 8).Value = 3757.6667  # H136: 4193.375 -> 3757.6667
$ws.Cells.Item(136, 9).Value = 3757.6667  # I136: 4309.4 -> 3757.6667
$ws.Cells.Item(136, 10).Value = 0  # J136: 4000 -> 0
$ws.Cells.Item(136, 11).Value = 11273.0001  # K136: 12928.2 -> 11273.0001
$ws.Cells.Item(136, 12).Value = 0  # L136: 12000 -> 0
$ws.Cells.Item(136, 13).Value = -8723.000100000001  # M136: -10378.2 -> -8723.000100000001
$ws.Cells.Item(136, 14).Value = ""  # N136: -17100 -> (blank)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 585.3570999999999  # H80: 555.73334 -> 585.3570999999999
$ws.Cells.Item(80, 9).Value = 310  # I80: 274 -> 310
$ws.Cells.Item(80, 10).Value = 952.5  # J80: 1119.2 -> 952.5
$ws.Cells.Item(80, 11).Value = 310  # K80: 274 -> 310
$ws.Cells.Item(80, 12).Value = 952.5  # L80: 1119.2 -> 952.5
$ws.Cells.Item(80, 13).Value = 688  # M80: 724 -> 688
$ws.Cells.Item(80, 14).Value = -2948.5  # N80: -3115.2 -> -2948.5
$ws.Cells.Item(83, 8).Value = 585.3570999999999  # H83: 555.73334 -> 585.3570999999999
$ws.Cells.Item(83, 9).Value = 310  # I83: 274 -> 310
$ws.Cells.Item(83, 10).Value = 952.5  # J83: 1119.2 -> 952.5
$ws.Cells.Item(83, 11).Value = 1550  # K83: 1370 -> 1550
$ws.Cells.Item(83, 12).Value = 4762.5  # L83: 5596 -> 4762.5
$ws.Cells.Item(83, 13).Value = 3442  # M83: 3622 -> 3442
$ws.Cells.Item(83, 14).Value = -14746.5  # N83: -15580 -> -14746.5
$ws.Cells.Item(94, 8).Value = 550.875  # H94: 546.2857 -> 550.875
$ws.Cells.Item(94, 9).Value = 550.875  # I94: 546.2857 -> 550.875
$ws.Cells.Item(94, 11).Value = 550.875  # K94: 546.2857 -> 550.875
$ws.Cells.Item(94, 13).Value = -99.875  # M94: -95.28570000000002 -> -99.875
$ws.Cells.Item(99, 8).Value = 3160.8333  # H99: 3162.1667 -> 3160.8333
$ws.Cells.Item(99, 9).Value = 2323  # I99: 2325.6667 -> 2323
$ws.Cells.Item(99, 11).Value = 2323  # K99: 2325.6667 -> 2323
$ws.Cells.Item(99, 13).Value = -825  # M99: -827.6667000000002 -> -825
$ws.Cells.Item(105, 8).Value = 2010  # H105: 1812.5 -> 2010
$ws.Cells.Item(105, 9).Value = 2010  # I105: 1727.25 -> 2010
$ws.Cells.Item(105, 10).Value = 0  # J105: 1983 -> 0
$ws.Cells.Item(105, 11).Value = 2010  # K105: 1727.25 -> 2010
$ws.Cells.Item(105, 12).Value = 0  # L105: 1983 -> 0
$ws.Cells.Item(105, 13).Value = -263  # M105: 19.75 -> -263
$ws.Cells.Item(105, 14).Value = ""  # N105: -5477 -> (blank)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(2, 8).Value = 272.8  # H2: 299.5 -> 272.8
$ws.Cells.Item(2, 9).Value = 85  # I2: 4 -> 85
$ws.Cells.Item(2, 11).Value = 85  # K2: 4 -> 85
$ws.Cells.Item(2, 13).Value = 28  # M2: 109 -> 28
$ws.Cells.Item(4, 8).Value = 6003994  # H4: 7503742.5 -> 6003994
$ws.Cells.Item(4, 10).Value = 5000  # J4: 0 -> 5000
$ws.Cells.Item(4, 12).Value = 5000  # L4: 0 -> 5000
$ws.Cells.Item(4, 14).Value = -5224  # N4: None -> -5224
$ws.Cells.Item(19, 8).Value = 136.8  # H19: 4000160 -> 136.8
$ws.Cells.Item(19, 9).Value = 160.25  # I19: 4000160 -> 160.25
$ws.Cells.Item(19, 10).Value = 43  # J19: 0 -> 43
$ws.Cells.Item(19, 11).Value = 160.25  # K19: 4000160 -> 160.25
$ws.Cells.Item(19, 12).Value = 43  # L19: 0 -> 43
$ws.Cells.Item(19, 13).Value = 9.75  # M19: -3999990 -> 9.75
$ws.Cells.Item(19, 14).Value = -383  # N19: None -> -383
$ws.Cells.Item(24, 8).Value = 136.8  # H24: 4000160 -> 136.8
$ws.Cells.Item(24, 9).Value = 160.25  # I24: 4000160 -> 160.25
$ws.Cells.Item(24, 10).Value = 43  # J24: 0 -> 43
$ws.Cells.Item(24, 11).Value = 160.25  # K24: 4000160 -> 160.25
$ws.Cells.Item(24, 12).Value = 43  # L24: 0 -> 43
$ws.Cells.Item(24, 13).Value = 9.75  # M24: -3999990 -> 9.75
$ws.Cells.Item(24, 14).Value = -383  # N24: None -> -383
$ws.Cells.Item(141, 8).Value = 69326  # H141: 70550.664 -> 69326
$ws.Cells.Item(141, 10).Value = 69326  # J141: 70550.664 -> 69326
$ws.Cells.Item(141, 12).Value = 69326  # L141: 70550.664 -> 69326
$ws.Cells.Item(141, 14).Value = -79686  # N141: -80910.664 -> -79686
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(54, 8).Value = 1150  # H54: 2649.4546 -> 1150
$ws.Cells.Item(54, 9).Value = 1300  # I54: 2814.4 -> 1300
$ws.Cells.Item(54, 11).Value = 3900  # K54: 8443.200000000001 -> 3900
$ws.Cells.Item(54, 13).Value = -3341  # M54: -7884.200000000001 -> -3341
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 8500  # H113: 8000 -> 8500
$ws.Cells.Item(113, 9).Value = 7000  # I113: 0 -> 7000
$ws.Cells.Item(113, 10).Value = 8800  # J113: 8000 -> 8800
$ws.Cells.Item(113, 11).Value = 7000  # K113: 0 -> 7000
$ws.Cells.Item(113, 12).Value = 8800  # L113: 8000 -> 8800
$ws.Cells.Item(113, 13).Value = -4830  # M113: None -> -4830
$ws.Cells.Item(113, 14).Value = -13140  # N113: -12340 -> -13140
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 11555.444  # H2: 11749.5 -> 11555.444
$ws.Cells.Item(2, 10).Value = 11555.444  # J2: 11749.5 -> 11555.444
$ws.Cells.Item(2, 12).Value = 11555.444  # L2: 11749.5 -> 11555.444
$ws.Cells.Item(2, 14).Value = -11779.444  # N2: -11973.5 -> -11779.444
$ws.Cells.Item(11, 8).Value = 0  # H11: 20000 -> 0
$ws.Cells.Item(11, 10).Value = 0  # J11: 20000 -> 0
$ws.Cells.Item(11, 12).Value = 0  # L11: 20000 -> 0
$ws.Cells.Item(11, 14).Value = ""  # N11: -20280 -> (blank)
$ws.Cells.Item(16, 8).Value = 1982.8334  # H16: 2500.5 -> 1982.8334
$ws.Cells.Item(16, 9).Value = 2079.4  # I16: 2500.5 -> 2079.4
$ws.Cells.Item(16, 10).Value = 1500  # J16: 0 -> 1500
$ws.Cells.Item(16, 11).Value = 2079.4  # K16: 2500.5 -> 2079.4
$ws.Cells.Item(16, 12).Value = 1500  # L16: 0 -> 1500
$ws.Cells.Item(16, 13).Value = -1909.4  # M16: -2330.5 -> -1909.4
$ws.Cells.Item(16, 14).Value = -1840  # N16: None -> -1840
$ws.Cells.Item(22, 8).Value = 2292.8572  # H22: 2312.5 -> 2292.8572
$ws.Cells.Item(22, 10).Value = 1790  # J22: 1900 -> 1790
$ws.Cells.Item(22, 12).Value = 1790  # L22: 1900 -> 1790
$ws.Cells.Item(22, 14).Value = -2380  # N22: -2490 -> -2380
$ws.Cells.Item(27, 8).Value = 2292.8572  # H27: 2312.5 -> 2292.8572
$ws.Cells.Item(27, 10).Value = 1790  # J27: 1900 -> 1790
$ws.Cells.Item(27, 12).Value = 1790  # L27: 1900 -> 1790
$ws.Cells.Item(27, 14).Value = -2004  # N27: -2114 -> -2004
$ws.Cells.Item(30, 8).Value = 474.75  # H30: 500 -> 474.75
$ws.Cells.Item(30, 9).Value = 474.75  # I30: 500 -> 474.75
$ws.Cells.Item(30, 11).Value = 474.75  # K30: 500 -> 474.75
$ws.Cells.Item(30, 13).Value = -366.75  # M30: -392 -> -366.75
$ws.Cells.Item(68, 8).Value = 5056.5  # H68: 5336.5454 -> 5056.5
$ws.Cells.Item(68, 9).Value = 1988  # I68: 2067.3333 -> 1988
$ws.Cells.Item(68, 10).Value = 8125  # J68: 6562.5 -> 8125
$ws.Cells.Item(68, 11).Value = 1988  # K68: 2067.3333 -> 1988
$ws.Cells.Item(68, 12).Value = 8125  # L68: 6562.5 -> 8125
$ws.Cells.Item(68, 13).Value = -1239  # M68: -1318.3333 -> -1239
$ws.Cells.Item(68, 14).Value = -9623  # N68: -8060.5 -> -9623
$ws.Cells.Item(71, 8).Value = 5056.5  # H71: 5336.5454 -> 5056.5
$ws.Cells.Item(71, 9).Value = 1988  # I71: 2067.3333 -> 1988
$ws.Cells.Item(71, 10).Value = 8125  # J71: 6562.5 -> 8125
$ws.Cells.Item(71, 11).Value = 9940  # K71: 10336.6665 -> 9940
$ws.Cells.Item(71, 12).Value = 40625  # L71: 32812.5 -> 40625
$ws.Cells.Item(71, 13).Value = -6196  # M71: -6592.666499999999 -> -6196
$ws.Cells.Item(71, 14).Value = -48113  # N71: -40300.5 -> -48113
$ws.Cells.Item(93, 8).Value = 1300.1  # H93: 1389 -> 1300.1
$ws.Cells.Item(93, 9).Value = 1429.1428  # I93: 1584 -> 1429.1428
$ws.Cells.Item(93, 11).Value = 1429.1428  # K93: 1584 -> 1429.1428
$ws.Cells.Item(93, 13).Value = -181.1428000000001  # M93: -336 -> -181.1428000000001
$ws.Cells.Item(100, 8).Value = 5352.5  # H100: 5352.5625 -> 5352.5
$ws.Cells.Item(100, 10).Value = 8855.333000000001  # J100: 8855.444 -> 8855.333000000001
$ws.Cells.Item(100, 12).Value = 8855.333000000001  # L100: 8855.444 -> 8855.333000000001
$ws.Cells.Item(100, 14).Value = -9937.333000000001  # N100: -9937.444 -> -9937.333000000001
$ws.Cells.Item(132, 8).Value = 2233.5557  # H132: 2585.6667 -> 2233.5557
$ws.Cells.Item(132, 9).Value = 1655.2858  # I132: 1749.75 -> 1655.2858
$ws.Cells.Item(132, 11).Value = 4965.857400000001  # K132: 5249.25 -> 4965.857400000001
$ws.Cells.Item(132, 13).Value = -2435.857400000001  # M132: -2719.25 -> -2435.857400000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1450  # H100: 1540 -> 1450
$ws.Cells.Item(100, 10).Value = 1833.3334  # J100: 2250 -> 1833.3334
$ws.Cells.Item(100, 12).Value = 3666.6668  # L100: 4500 -> 3666.6668
$ws.Cells.Item(100, 14).Value = -4748.6668  # N100: -5582 -> -4748.6668
